$wb = $excel.ActiveWorkbook

# --- Update the localization status text from "Ready for handoff" to
#     "In Translation" everywhere it appears (Overview!E2:F2, and the
#     "Status" column on the per-language sheets).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
